$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3488
$ws.Range("L3").Value = 3637
$ws.Range("I4").Value = 1840
$ws.Range("J4").Value = 1869
$ws.Range("L4").Value = 908
$ws.Range("L5").Value = 219
$ws.Range("L6").Value = 3186
$ws.Range("I7").Value = 26308
$ws.Range("J7").Value = 29345
$ws.Range("L7").Value = 11438

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L3").Value = 3
$ws.Range("L6").Value = 10

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 246
$ws.Range("L7").Value = 741

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 169
$ws.Range("L7").Value = 537

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 124
$ws.Range("L7").Value = 407

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 44
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 385
$ws.Range("L8").Value = 741
$ws.Range("L10").Value = 71
$ws.Range("L11").Value = 191
$ws.Range("L17").Value = 20
$ws.Range("L19").Value = 320
$ws.Range("L21").Value = 34
$ws.Range("L23").Value = 125
$ws.Range("L29").Value = 616
$ws.Range("L31").Value = 112
$ws.Range("L33").Value = 537
$ws.Range("L34").Value = 71
$ws.Range("L36").Value = 156
$ws.Range("L37").Value = 407
$ws.Range("L38").Value = 10
$ws.Range("L40").Value = 32
$ws.Range("L43").Value = 84
$ws.Range("L48").Value = 160
$ws.Range("L50").Value = 55
$ws.Range("L51").Value = 143
$ws.Range("L52").Value = 237
$ws.Range("L54").Value = 243
$ws.Range("L55").Value = 109
$ws.Range("I63").Value = 264
$ws.Range("J63").Value = 221
$ws.Range("L63").Value = 34
$ws.Range("L66").Value = 29
$ws.Range("L67").Value = 410
$ws.Range("L76").Value = 166
$ws.Range("L77").Value = 68
$ws.Range("L83").Value = 264
$ws.Range("L85").Value = 581
$ws.Range("L90").Value = 109
$ws.Range("L94").Value = 137
$ws.Range("L96").Value = 113
$ws.Range("L97").Value = 99
$ws.Range("L99").Value = 191
$ws.Range("I101").Value = 26308
$ws.Range("J101").Value = 29345
$ws.Range("L101").Value = 11438

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 119
$ws.Range("L3").Value = 154
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 410

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 58
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 239
$ws.Range("L4").Value = 33
$ws.Range("L7").Value = 616

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 109
$ws.Range("L7").Value = 320

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 31
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 47
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 156

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 108
$ws.Range("L7").Value = 385

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 72
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 32
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L2").Value = 18
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 168
$ws.Range("L3").Value = 239
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 581

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 82
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 237
